$p = $ppt.ActivePresentation

# --- Update the "last saved" datetimeFigureOut field text (master + every
#     slide layout) from 2015/12/4 to 2016/1/1. The date placeholder is the
#     shape whose PlaceholderFormat.Type equals 16 (ppPlaceholderDate).
$sm = $p.Designs.Item(1).SlideMaster

$msh = $sm.Shapes
for ($i = 1; $i -le $msh.Count; $i++) {
    $ms = $msh.Item($i)
    if ($ms.PlaceholderFormat.Type -eq 16) {
        $ms.TextFrame.TextRange.Text = "2016/1/1"
    }
}

$layouts = $sm.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    $lay = $layouts.Item($j)
    $lsh = $lay.Shapes
    for ($i = 1; $i -le $lsh.Count; $i++) {
        $ls = $lsh.Item($i)
        if ($ls.PlaceholderFormat.Type -eq 16) {
            $ls.TextFrame.TextRange.Text = "2016/1/1"
        }
    }
}

# --- Slide 2: widen the "Filter critical: ..." textbox and extend its
#     wording with ", covered" after "identity".
$s2 = $p.Slides.Item(2)
$sh = $s2.Shapes.Item(2)
$sh.Width = 428.9484251968504

$tr = $sh.TextFrame.TextRange
$tail = $tr.Characters(35, 17)
$tail.Text = ", gaps, identity, covered "
$idRun = $tr.Characters(43, 8)
$idRun.Text = "identity"
$spRun = $tr.Characters(60, 1)
$spRun.Text = " "

Write-Host "done"
